$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F; this shifts old F(Grade)->G, old G(Weighted formula)->H
$ws.Columns.Item(6).Insert()

# New column headers
$ws.Range("F3").Value = "Weight of CWK"
$ws.Range("I3").Value = "CWK"

# Set "Weight of CWK" values (column F) for rows that have a weight
$ws.Range("F8").Value = 0.625
$ws.Range("F12").Value = 0.625
$ws.Range("F14").Value = 0.625
$ws.Range("F15").Value = 0.625
$ws.Range("F16").Value = 0.625
$ws.Range("F17").Value = 12.5
$ws.Range("F18").Value = 0.625
$ws.Range("F19").Value = 0.625
$ws.Range("F20").Value = 0.625
$ws.Range("F21").Value = 0.625
$ws.Range("F23").Value = 0.625
$ws.Range("F24").Value = 75
$ws.Range("F25").Value = 6.25

# Row 17: Grade (G17) gets a new value of 69
$ws.Range("G17").Value = 69

# Row 13: clear E13 and the shifted-formula cell H13; keep only G13 = 64
$ws.Range("E13").ClearContents()
$ws.Range("H13").ClearContents()

# Row 22: clear E22 and the shifted-formula cell H22 entirely
$ws.Range("E22").ClearContents()
$ws.Range("H22").ClearContents()

# Add CWK formulas in column I for the applicable rows
$ws.Range("I8").Formula = "=F8*G8/100"
$ws.Range("I12").Formula = "=F12*G12/100"
$ws.Range("I14").Formula = "=F14*G14/100"
$ws.Range("I15").Formula = "=F15*G15/100"
$ws.Range("I16").Formula = "=F16*G16/100"
$ws.Range("I17").Formula = "=F17*G17/100"
$ws.Range("I18").Formula = "=F18*G18/100"
$ws.Range("I19").Formula = "=F19*G19/100"
$ws.Range("I20").Formula = "=F20*G20/100"
$ws.Range("I21").Formula = "=F21*G21/100"
$ws.Range("I23").Formula = "=F23*G23/100"
$ws.Range("I24").Formula = "=F24*G24/100"
$ws.Range("I25").Formula = "=F25*G25/100"

# Row 27 sum formulas
$ws.Range("H27").Formula = "=SUM(H4:H26)"
$ws.Range("I27").Formula = "=SUM(I8:I26)"

# Column F width (approx bestFit)
$ws.Columns.Item(6).ColumnWidth = 12.43

# Update selection to match final cursor position
$ws.Range("N24").Select()
